# Generate Report for Archive
# Swap the report rows for ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.md and
# 308ccd35-ecfe-48ef-a69d-f568ccae4764.md on all three sheets.

$wb = $excel.ActiveWorkbook

function Swap-HyperlinkText($ws, $oldText, $newText) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.TextToDisplay -eq $oldText) {
            $h.TextToDisplay = $newText
        }
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = "308ccd35-ecfe-48ef-a69d-f568ccae4764.md"
$wsOverview.Range("A4").Value = "ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.md"
$wsOverview.Range("B4").Value = "In Translation"
$wsOverview.Range("C4").Value = "In Translation"

Swap-HyperlinkText $wsOverview "ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.md" "308ccd35-ecfe-48ef-a69d-f568ccae4764.md"
Swap-HyperlinkText $wsOverview "308ccd35-ecfe-48ef-a69d-f568ccae4764.md" "ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.md"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value = "308ccd35-ecfe-48ef-a69d-f568ccae4764.md"
$wsZh.Range("C3").Value = "308ccd35-ecfe-48ef-a69d-f568ccae4764.ea0ffc17f3e214385cd419df89c6c1e8b2d7a6df.zh-cn.xlf"
$wsZh.Range("D3").Value = "2016-02-24 09:27:40"

$wsZh.Range("A4").Value = "ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.md"
$wsZh.Range("B4").Value = "In Translation"
$wsZh.Range("C4").Value = "ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.3e14b8dd899da539dd170c68c834efb820e8f44c.zh-cn.xlf"
$wsZh.Range("D4").Value = "2016-02-24 09:25:45"

Swap-HyperlinkText $wsZh "ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.md" "308ccd35-ecfe-48ef-a69d-f568ccae4764.md"
Swap-HyperlinkText $wsZh "308ccd35-ecfe-48ef-a69d-f568ccae4764.md" "ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.md"
Swap-HyperlinkText $wsZh "ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.3e14b8dd899da539dd170c68c834efb820e8f44c.zh-cn.xlf" "308ccd35-ecfe-48ef-a69d-f568ccae4764.ea0ffc17f3e214385cd419df89c6c1e8b2d7a6df.zh-cn.xlf"
Swap-HyperlinkText $wsZh "308ccd35-ecfe-48ef-a69d-f568ccae4764.ea0ffc17f3e214385cd419df89c6c1e8b2d7a6df.zh-cn.xlf" "ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.3e14b8dd899da539dd170c68c834efb820e8f44c.zh-cn.xlf"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value = "308ccd35-ecfe-48ef-a69d-f568ccae4764.md"
$wsDe.Range("C3").Value = "308ccd35-ecfe-48ef-a69d-f568ccae4764.ea0ffc17f3e214385cd419df89c6c1e8b2d7a6df.de-de.xlf"
$wsDe.Range("D3").Value = "2016-02-24 09:27:52"

$wsDe.Range("A4").Value = "ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.md"
$wsDe.Range("B4").Value = "In Translation"
$wsDe.Range("C4").Value = "ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.3e14b8dd899da539dd170c68c834efb820e8f44c.de-de.xlf"
$wsDe.Range("D4").Value = "2016-02-24 09:25:57"

Swap-HyperlinkText $wsDe "ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.md" "308ccd35-ecfe-48ef-a69d-f568ccae4764.md"
Swap-HyperlinkText $wsDe "308ccd35-ecfe-48ef-a69d-f568ccae4764.md" "ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.md"
Swap-HyperlinkText $wsDe "ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.3e14b8dd899da539dd170c68c834efb820e8f44c.de-de.xlf" "308ccd35-ecfe-48ef-a69d-f568ccae4764.ea0ffc17f3e214385cd419df89c6c1e8b2d7a6df.de-de.xlf"
Swap-HyperlinkText $wsDe "308ccd35-ecfe-48ef-a69d-f568ccae4764.ea0ffc17f3e214385cd419df89c6c1e8b2d7a6df.de-de.xlf" "ff7f135a-bb2f-4e7b-8b06-4236cca36a5a.3e14b8dd899da539dd170c68c834efb820e8f44c.de-de.xlf"
